# Apply the data updates described by the diff.
# Columns A, B, E, F, G, H, Q, R in rows 2, 3, 5, 6, 7 are updated to new values
# (the rows were effectively cyclically rotated: old row 6 -> row 2, old row 7 -> row 3,
#  old row 2 -> row 5, old row 3 -> row 6, old row 5 -> row 7).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 111525238
$ws.Range("B2").Value = 77515
$ws.Range("E2").Value = 6425
$ws.Range("F2").Value = "Garnlav"
$ws.Range("G2").Value = "Alectoria sarmentosa"
$ws.Range("H2").Value = "(Ach.) Ach."
$ws.Range("Q2").Value = 404495.4563026094
$ws.Range("R2").Value = 6706677.491168984

$ws.Range("A3").Value = 111525224
$ws.Range("B3").Value = 77268
$ws.Range("E3").Value = 228912
$ws.Range("F3").Value = "Mörk kolflarnlav"
$ws.Range("G3").Value = "Carbonicola myrmecina"
$ws.Range("H3").Value = "(Ach.) Bendiksby & Timdal"
$ws.Range("Q3").Value = 404619.9854206198
$ws.Range("R3").Value = 6706773.322858612

$ws.Range("A5").Value = 111525223
$ws.Range("B5").Value = 78107
$ws.Range("E5").Value = 6453
$ws.Range("F5").Value = "Vedskivlav"
$ws.Range("G5").Value = "Hertelidea botryosa"
$ws.Range("H5").Value = "(Fr.) Printzen & Kantvilas"
$ws.Range("Q5").Value = 404637.0659126193
$ws.Range("R5").Value = 6706784.214121711

$ws.Range("A6").Value = 111525235
$ws.Range("B6").Value = 77515
$ws.Range("E6").Value = 6425
$ws.Range("F6").Value = "Garnlav"
$ws.Range("G6").Value = "Alectoria sarmentosa"
$ws.Range("H6").Value = "(Ach.) Ach."
$ws.Range("Q6").Value = 404485.2245768273
$ws.Range("R6").Value = 6706757.647421388

$ws.Range("A7").Value = 111525226
$ws.Range("B7").Value = 77515
$ws.Range("E7").Value = 6425
$ws.Range("F7").Value = "Garnlav"
$ws.Range("G7").Value = "Alectoria sarmentosa"
$ws.Range("H7").Value = "(Ach.) Ach."
$ws.Range("Q7").Value = 404616.9589749529
$ws.Range("R7").Value = 6706770.937089294
